$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Landesliga")

# The data refresh swapped the recorded results/odds between the match shown on
# row 11 (id 9) and the match shown on row 13 (id 11), while each row's "id" and
# "HomeTeam" stayed put. Concretely: the match id, AwayTeam, score and all odds
# columns (H through AC, except R/S/X/AB which are identical on both rows) were
# exchanged between the two rows.

# --- New values for row 11 (take former row 13 values) ---
$ws.Range("B11").Value = 7035048
$ws.Range("G11").Value = "TuRU Dsseldorf"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "H"
$ws.Range("K11").Value = 3.25
$ws.Range("L11").Value = 4
$ws.Range("M11").Value = 1.8
$ws.Range("N11").Value = 2.9
$ws.Range("O11").Value = 4
$ws.Range("P11").Value = 1.95
$ws.Range("Q11").Value = 0.5
$ws.Range("T11").Value = 3
$ws.Range("U11").Value = 1.75
$ws.Range("V11").Value = 1.95
$ws.Range("W11").Value = 1.9
$ws.Range("Y11").Value = -1
$ws.Range("Z11").Value = 0.8
$ws.Range("AA11").Value = -1
$ws.Range("AC11").Value = 0.95

# --- New values for row 13 (take former row 11 values) ---
$ws.Range("B13").Value = 7035046
$ws.Range("G13").Value = "FC Viersen"
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = "A"
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 3.6
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 2
$ws.Range("O13").Value = 3.6
$ws.Range("P13").Value = 3
$ws.Range("Q13").Value = -0.25
$ws.Range("T13").Value = 2.75
$ws.Range("U13").Value = 1.8
$ws.Range("V13").Value = 2
$ws.Range("W13").Value = -1
$ws.Range("Y13").Value = 2
$ws.Range("Z13").Value = -1
$ws.Range("AA13").Value = 1
$ws.Range("AC13").Value = 1

$wb.Save()
